# ---------------------------------------------------------------------------
# Applies the "Balancing Life and Technology" -> "Unraveling the Art of
# Chemistry" rewrite described by the supplied OOXML diff.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Replace-One([string]$old, [string]$new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $ok) {
        Write-Host "WARNING: could not find text: $old"
    }
    return $ok
}

# Inserts $extra (which should start with the punctuation/space that
# continues the sentence, e.g. ". Next sentence") right after the text
# identified by $anchor (the *new*, already-replaced, sentence text) and
# before whatever immediately follows it (typically the pre-existing
# "." run) so the existing trailing punctuation ends up after our new text.
# NOTE: deliberately does NOT touch $ins.Font.* - writing a Font property
# on a zero-length (collapsed) Range mis-targets unrelated runs in this
# runtime, so we just let InsertAfter inherit the surrounding formatting.
function Insert-After([string]$anchor, [string]$extra) {
    $r = $d.Content
    $ok = $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "WARNING: could not find anchor text: $anchor"
        return
    }
    $pos = $r.End
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($extra)
}

# ---------------------------------------------------------------------------
# Title / author / email
# ---------------------------------------------------------------------------

Replace-One "Balancing Life and Technology" "Unraveling the Art of Chemistry: A Story of Elements and Transformation"
Replace-One "Stefanie Campbell" "Alexis Robertson"
Replace-One "adrienne@aenable" "alexis"
Replace-One "net" "r@educonnect"
Insert-After "r@educonnect" ".org"

# ---------------------------------------------------------------------------
# First body paragraph (essay body)
# ---------------------------------------------------------------------------

Replace-One "In our increasingly digital world, technology has become an integral part of our lives" "Chemistry is the science that explores the composition, structure, properties, and change of matter"

Replace-One " We rely on it for communication, education, and entertainment" " It is a fundamental discipline that touches various aspects of life, from the food we eat and the clothes we wear to the medicines we take and the technologies we use"

Replace-One " While technology brings many benefits, it also raises concerns about its potential to disrupt our lives and disconnect us from the real world" " Chemistry reveals the fascinating world of atoms, molecules, and their interactions, unlocking secrets hidden within the material universe"

Replace-One " In this essay, we explore the delicate balance between embracing technological advancements and maintaining a healthy, fulfilling life" " By studying chemistry, we gain an understanding of how these tiny particles combine to form substances with diverse properties, leading to transformative applications across industries and domains"
Insert-After " By studying chemistry, we gain an understanding of how these tiny particles combine to form substances with diverse properties, leading to transformative applications across industries and domains" ". From the birth of stars to the chemical reactions within our bodies, chemistry paints a vivid tapestry of the universe's intricate workings"

Replace-One "We begin by examining the ways in which technology has transformed our lives" "Furthermore, chemistry has played a pivotal role in shaping human history"

Replace-One " From smartphones to social media to artificial intelligence, technology has changed the way we learn, connect, work, and play" " Early civilizations relied on chemical knowledge to craft tools, preserve food, and develop medicines, gradually evolving into the advanced chemical industries of today"

Replace-One " We discuss both the positive and negative impacts of these changes, highlighting the importance of using technology mindfully and responsibly" " From the discovery of elements to the synthesis of complex molecules, chemistry has revolutionized the way we live, work, and interact with the environment"
Insert-After " From the discovery of elements to the synthesis of complex molecules, chemistry has revolutionized the way we live, work, and interact with the environment" ". It has enabled the creation of groundbreaking materials, fertilizers that feed nations, and pharmaceuticals that combat diseases, profoundly impacting global health and well-being. As we continue to explore the depths of chemical processes, we delve into mysteries yet unsolved, unlocking the potential for transformative discoveries that will continue to shape the future of science and technology"

Replace-One "Next, we delve into the psychological and social effects of excessive technology use" "In delving into the vast realm of chemistry, we embark on a journey filled with intrigue and revelation"

Replace-One " We explore the concept of tech addiction, discussing the addictive nature of certain technologies and the negative consequences they can have on our mental and emotional well-being" " We witness the magic of atoms rearranging to form new substances, the dance of molecules as they undergo reactions, and the symphony of energy transformations that drive chemical processes"

Replace-One " We also examine the impact of technology on our relationships, both online and offline, and how it can lead to isolation and decreased face-to-face interactions" " With each discovery, we unveil a piece of the intricate puzzle that governs the behavior of matter, gaining insights into the essence of the natural world and our place within it"
Insert-After " With each discovery, we unveil a piece of the intricate puzzle that governs the behavior of matter, gaining insights into the essence of the natural world and our place within it" ". Whether unraveling the mysteries of life through biochemistry or unlocking the secrets of material properties, chemistry invites us to explore the boundless wonders that lie at the heart of our universe"

# ---------------------------------------------------------------------------
# Summary heading / body
# ---------------------------------------------------------------------------

Replace-One "In this essay, we examined the complex relationship between life and technology" "Chemistry is a fascinating science that explores the composition, structure, properties, and change of matter"

Replace-One " We explored the ways in which technology has transformed our lives, highlighting both the benefits and the risks associated with its use" " It has played a crucial role in shaping human history, leading to the development of tools, medicines, and advanced technologies"

Replace-One " We delved into the psychological and social effects of excessive technology use, discussing the concept of tech addiction and its impact on our mental and emotional well-being" " Chemistry has revolutionized industries, transformed global health, and continues to unveil the mysteries of life and the material world"

Replace-One " By understanding the potential consequences of technology overuse, we can take steps to find a balance that allows us to enjoy its benefits without sacrificing our health and relationships" " By studying chemistry, we gain an understanding of the fundamental principles that govern the universe, unlocking the potential for transformative discoveries that will shape the future of science and technology"

# ---------------------------------------------------------------------------
# Trailing empty paragraph added at the end of the document body
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
